# Fix Plots.xlsx data names
#
# The "DataCombined" sheet stores a dataset name that was missing an
# underscore between "Human" and "PeripheralVenousBlood". Correct the
# two cells (F3 and F5) that hold this value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataCombined")

$oldValue = "Laskin 1982.Group A_Aciclovir_1_Human_PeripheralVenousBlood_Plasma_2.5 mg/kg_iv_"
$newValue = "Laskin 1982.Group A_Aciclovir_1_Human__PeripheralVenousBlood_Plasma_2.5 mg/kg_iv_"

foreach ($cellAddr in @("F3", "F5")) {
    $cell = $ws.Range($cellAddr)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
